# Update faturamento_diario - "atualizei dados da bibi e add"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing values for days 8 and 10 (July/2025)
$ws.Range("B7").Value = 20626.55
$ws.Range("B9").Value = 21120.89

# Insert a new row for day 11 (July/2025) before the June data block,
# shifting all subsequent rows down by one.
$ws.Rows("10:10").Insert()

$ws.Range("A10").Value = 11
$ws.Range("B10").Value = 13497.7
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 2025
$ws.Range("E10").Value = "07/2025"
